$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new "description" header in column M, row 1 (new last column after "order")
$ws.Range("M1").Value = "description"

# Update selection to just M1 (was previously the whole column M1:M1048576)
$ws.Range("M1").Select()
